$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 13 only held the "4808662 - Lucrécio Fábio dos Santos" value in B/C
# (column A was blank there). Removing it shifts every following row up by one,
# which is exactly what the target layout (C1:C24 instead of C1:C25) needs.
$ws.Rows.Item(13).Delete()

# After the shift, a handful of B/C cells need their text replaced to match the
# new content laid out in the edited workbook.
$ws.Range("B10").Value = "4808662 - Lucrécio Fábio dos Santos"
$ws.Range("C10").Value = "4808662 - Lucrécio Fábio dos Santos"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Value = "01/01/2021"
$ws.Range("C15").Value = "01/01/2021"

$ws.Range("B18").Value = "4808662 - Lucrécio Fábio dos Santos"
$ws.Range("C18").Value = "4808662 - Lucrécio Fábio dos Santos"

$ws.Range("B19").Value = "Aulas expositivas teóricas, aulas práticas, aulas de exercícios, aulas de laboratório."
$ws.Range("C19").Value = "Aulas expositivas teóricas, aulas práticas, aulas de exercícios, aulas de laboratório."

$ws.Range("B20").Value = "Nota de duas provas (P1 e P2)Fórmula: M1 = (P1 + 2 x P2)/3.."
$ws.Range("C20").Value = "Nota de duas provas (P1 e P2)Fórmula: M1 = (P1 + 2 x P2)/3.."

$ws.Range("B21").Value = "Aplicação de uma prova envolvendo o assunto de todo semestre.NR (nota da recuperação) = (M1 + NR)/2."
$ws.Range("C21").Value = "Aplicação de uma prova envolvendo o assunto de todo semestre.NR (nota da recuperação) = (M1 + NR)/2."
